# Updated cryptos list on Mon Jun 10 21:48:53 UTC 2024 with GitHub Actions
# Refreshes the Price/Volume(1h) figures pulled from coinranking.com and
# reorders the Kaspa / NEARProtocol rows to match the latest ranking.
#
# Note: several "Price" values look like plain decimal numbers (e.g.
# "621.52"). The source data stores these as plain text (no thousands
# separators use a dot+dot pattern like "69.598.78" which is naturally
# text, but single-dot values parse as numbers). To keep them as text -
# exactly like the original workbook - NumberFormat is forced to "@"
# before the assignment, then the style is reset back to Normal so no
# extra formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.598.78"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.672.63"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "621.52"
$ws.Range("E5").Value = "  -7.39%  "
Set-TextValue "D6" "159.21"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue "D8" "0.496"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  +1.34%  "
Set-TextValue "D11" "0.439"
$ws.Range("E11").Value = "  -1.15%  "
Set-TextValue "D12" "0.0000229"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").Value = "4.290.55"
$ws.Range("E13").Value = "  -0.74%  "
Set-TextValue "D14" "32.36"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "3.652.62"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "69.621.15"
$ws.Range("E17").Value = "  +0.43%  "
Set-TextValue "D18" "6.50"
$ws.Range("E18").Value = "  -0.01%  "
Set-TextValue "D19" "15.89"
$ws.Range("E19").Value = "  -2.31%  "
Set-TextValue "D20" "10.32"
$ws.Range("E20").Value = "  +5.23%  "
Set-TextValue "D21" "470.30"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("E22").Value = "  -0.81%  "
Set-TextValue "D23" "79.79"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "3.817.41"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -3.26%  "
Set-TextValue "D27" "11.04"
$ws.Range("E27").Value = "  +0.90%  "
Set-TextValue "D28" "8.72"
$ws.Range("E28").Value = "  -4.25%  "
Set-TextValue "D29" "2.60"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("E30").Value = "  -4.44%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -1.19%  "

# Rows 34/35 swapped places (Kaspa <-> NEARProtocol) with refreshed figures
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "6.39"
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D35" "0.163"
$ws.Range("E35").Value = "  -2.52%  "

$ws.Range("D36").Value = "3.672.34"
$ws.Range("E36").Value = "  -0.48%  "
Set-TextValue "D37" "8.27"
$ws.Range("E37").Value = "  -3.30%  "
Set-TextValue "D39" "177.83"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -4.98%  "
Set-TextValue "D42" "2.21"
$ws.Range("E42").Value = "  -1.88%  "
Set-TextValue "D43" "0.0891"
$ws.Range("E43").Value = "  -2.13%  "
Set-TextValue "D44" "0.925"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  -0.80%  "
Set-TextValue "D46" "28.84"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("E47").Value = "  -2.63%  "
Set-TextValue "D48" "7.85"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E49").Value = "  -6.61%  "
$ws.Range("E50").Value = "  -4.79%  "
Set-TextValue "D51" "1.21"
$ws.Range("E51").Value = "  -5.89%  "
